$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in James Morse's (row 8) weekly availability, which was previously blank.
$ws.Range("B8").Value = "9am-1pm AND 2pm-MN"
$ws.Range("C8").Value = "1pm-MN"
$ws.Range("D8").Value = "9am-1pm AND 2pm-MN"
$ws.Range("E8").Value = "1pm-MN"
$ws.Range("F8").Value = "9am-1pm AND 2pm-MN"
$ws.Range("G8").Value = "8am-MN"
$ws.Range("H8").Value = "8am-MN"

# Move the selection to reflect where the user ended up after entering data.
$ws.Range("H8").Select()
